# "Actualizamos el minuto en la rama adicional"
# The paragraph "Continuar en el minuto 14:58 creando ramas" needs its
# timestamp updated from 14:58 to 18:50. In the target OOXML, the run
# that used to hold "14:58 creando ramas" is split in two: "18:50 " and
# "creando ramas" (same character formatting: red text, es-ES language).

$d = $word.ActiveDocument

# 1) Update the timestamp text itself.
$d.Content.Find.Execute("14:58", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "18:50", 2) | Out-Null

# 2) Force a run boundary right after "18:50 " (before "creando ramas")
#    by toggling a character property back to its original value; this
#    splits what would otherwise remain a single merged run into the
#    two runs the target document expects.
$tail = $d.Content.Duplicate
$tail.Find.Execute("creando ramas") | Out-Null
$tail.Font.Bold = 1
$tail.Font.Bold = 0

# 3) Likewise make sure "18:50 " is its own run, distinct from the
#    preceding "en el minuto " run.
$time = $d.Content.Duplicate
$time.Find.Execute("18:50 ") | Out-Null
$time.Font.Bold = 1
$time.Font.Bold = 0
